$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update shared-string templated cell values in row 2
$ws.Range("E2").Value = '${data.birthday}'
$ws.Range("O2").Value = '${data.studyOn}'
$ws.Range("P2").Value = '${data.graduateOn}'
$ws.Range("R2").Value = '${data.president}'

# 2. Change the date number format (yyyymmdd) to plain text (@) for the
#    birthday / studyOn / graduateOn column cells in row 2
$ws.Range("E2").NumberFormat = "@"
$ws.Range("O2").NumberFormat = "@"
$ws.Range("P2").NumberFormat = "@"

# 3. Add a new (empty) row 3 with the same text format on E3/O3/P3
$ws.Range("E3").NumberFormat = "@"
$ws.Range("O3").NumberFormat = "@"
$ws.Range("P3").NumberFormat = "@"

# 4. Bump the outline level so sheetFormatPr/outlineLevelRow becomes 2
$ws.Rows.Item(3).OutlineLevel = 2

# 5. Move/save the active selection to E3
$ws.Range("E3").Select()
